$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing GDP values (column B, rows 12-30)
$ws.Range("B12").Value = 107439.627
$ws.Range("B13").Value = 109575.857
$ws.Range("B14").Value = 111826.233
$ws.Range("B15").Value = 115122.352
$ws.Range("B16").Value = 123337.59
$ws.Range("B17").Value = 129242.845
$ws.Range("B18").Value = 138503.714
$ws.Range("B19").Value = 140775.522
$ws.Range("B20").Value = 138281.399
$ws.Range("B21").Value = 142517.164
$ws.Range("B22").Value = 147392.855
$ws.Range("B23").Value = 154899.048
$ws.Range("B24").Value = 164035.645
$ws.Range("B25").Value = 173826.485
$ws.Range("B26").Value = 184611.894
$ws.Range("B27").Value = 190857.987
$ws.Range("B28").Value = 199940.107
$ws.Range("B29").Value = 213146.127
$ws.Range("B30").Value = 227259.597

# Add new row 31 with next year's data
$ws.Range("A31").Value = 43831
$ws.Range("A31").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B31").Value = 223145.76
$ws.Range("B31").NumberFormat = "0.000"

# Update selection to match post-edit state (full A:B column selection)
$ws.Columns("A:B").Select()
